# 1) Slide 16 (index 16): table's style was changed to a different built-in
#    PowerPoint table style (GUID swap on <a:tableStyleId>).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
$tbl.ApplyStyle("{9AC5380E-C8BA-459F-95C8-1CABE32011E8}")

# 2) The deck's applied theme (ppt/theme/theme2.xml, used by the one Slide
#    Master) had its color scheme swapped from the custom "Integral" palette
#    to the stock Office palette (the "Office Theme"/"Integral" theme parts
#    were exchanged). Reproduce that by rewriting the twelve theme colors on
#    the live Slide Master's theme color scheme to the standard Office
#    theme values (font scheme / format scheme are already identical between
#    the two theme parts, so only the color scheme needs to change).
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# PowerPoint COM ColorFormat.RGB values are packed as 0x00BBGGRR (i.e. the
# same little-endian byte order VBA's RGB() function produces), so build
# each value from its R/G/B bytes explicitly (no RGB() helper in PowerShell).
function BGR([int]$r, [int]$g, [int]$b) { return $r + ($g * 256) + ($b * 65536) }

$cs.Colors(1).RGB  = BGR 0x00 0x00 0x00   # dk1
$cs.Colors(2).RGB  = BGR 0xFF 0xFF 0xFF   # lt1
$cs.Colors(3).RGB  = BGR 0x44 0x54 0x6A   # dk2
$cs.Colors(4).RGB  = BGR 0xE7 0xE6 0xE6   # lt2
$cs.Colors(5).RGB  = BGR 0x5B 0x9B 0xD5   # accent1
$cs.Colors(6).RGB  = BGR 0xED 0x7D 0x31   # accent2
$cs.Colors(7).RGB  = BGR 0xA5 0xA5 0xA5   # accent3
$cs.Colors(8).RGB  = BGR 0xFF 0xC0 0x00   # accent4
$cs.Colors(9).RGB  = BGR 0x44 0x72 0xC4   # accent5
$cs.Colors(10).RGB = BGR 0x70 0xAD 0x47   # accent6
$cs.Colors(11).RGB = BGR 0x05 0x63 0xC1   # hlink
$cs.Colors(12).RGB = BGR 0x95 0x4F 0x72   # folHlink
